$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2, 3 and 5 (identified by columns A, Q, R, AC) got
# cyclically rotated: the values that used to be on row 5 moved to
# row 2, the values on row 2 moved to row 3, and the values on row 3
# moved to row 5. Capture the original values first, then write the
# rotated values back.

$A2 = $ws.Range("A2").Value2
$Q2 = $ws.Range("Q2").Value2
$R2 = $ws.Range("R2").Value2
$AC2 = $ws.Range("AC2").Value2

$A3 = $ws.Range("A3").Value2
$Q3 = $ws.Range("Q3").Value2
$R3 = $ws.Range("R3").Value2
$AC3 = $ws.Range("AC3").Value2

$A5 = $ws.Range("A5").Value2
$Q5 = $ws.Range("Q5").Value2
$R5 = $ws.Range("R5").Value2
$AC5 = $ws.Range("AC5").Value2

# Row 2 gets the values that were previously on row 5
$ws.Range("A2").Value2 = $A5
$ws.Range("Q2").Value2 = $Q5
$ws.Range("R2").Value2 = $R5
$ws.Range("AC2").Value2 = $AC5

# Row 3 gets the values that were previously on row 2
$ws.Range("A3").Value2 = $A2
$ws.Range("Q3").Value2 = $Q2
$ws.Range("R3").Value2 = $R2
$ws.Range("AC3").Value2 = $AC2

# Row 5 gets the values that were previously on row 3
$ws.Range("A5").Value2 = $A3
$ws.Range("Q5").Value2 = $Q3
$ws.Range("R5").Value2 = $R3
$ws.Range("AC5").Value2 = $AC3
